# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Refresh the Albania MSME summary statistics with updated (more precise)
# source figures. All of these cells are stored as text in the workbook
# (they originated from a shared-strings text import), so we force each
# write to stay text-typed -- and restore the cell's original style
# afterwards so we don't leave a stray NumberFormat behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($rng, [string]$val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Source Type: Statistical Institution block -----------------------------

# Enterprises density (per 1000 people)
Set-TextValue $ws.Range("B11") "22.63"
Set-TextValue $ws.Range("C11") "2.14"
Set-TextValue $ws.Range("D11") "24.77"

# Employment (% of total)
Set-TextValue $ws.Range("B12") "36.82"
Set-TextValue $ws.Range("C12") "34.14"
Set-TextValue $ws.Range("D12") "70.96"

# Value added to the economy (% of total)
Set-TextValue $ws.Range("B18") "19.11"
Set-TextValue $ws.Range("C18") "40.19"

# Source Type: SME Associations (Most Widely Used) block -----------------

# Enterprises density (per 1000 people) -- D33 previously (erroneously)
# reused the "24.8" string from D11; it now gets its own distinct value.
Set-TextValue $ws.Range("B33") "23.76"
Set-TextValue $ws.Range("C33") "1.08"
Set-TextValue $ws.Range("D33") "24.84"

# Employment (% of total) -- B34 (44.9) is unchanged
Set-TextValue $ws.Range("C34") "36.47"
Set-TextValue $ws.Range("D34") "81.36"

# Enterprises (% of total)
Set-TextValue $ws.Range("B36") "95.54"
Set-TextValue $ws.Range("C36") "4.34"
Set-TextValue $ws.Range("D36") "99.89"
